# Updates cryptos price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.937.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.539.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.40%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.84%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.537.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.38%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").Value = '  -2.28%  '
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.993.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.891.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.547.89'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '333.43'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.83%  '
$ws.Range("E25").Value = '  -2.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.97%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  +11.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0813'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '177.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  +6.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '413.81'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.15%  '
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.83'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.17%  '
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '151.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.58'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.603'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.48%  '
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0517'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.86%  '
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("E51").Value = '  +2.60%  '
